$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Cell writes, ordered to reproduce the source workbook's shared-string table
# (new unique strings must be introduced in the same sequence as the original
# edit so sharedStrings.xml indices line up with the target diff).
# ---------------------------------------------------------------------------
$ws.Range("K1").Value = "nopol"
$ws.Range("J1").Value = "kontrak"
$ws.Range("J2").Value = "belum"
$ws.Range("G2").Value = 261294
$ws.Range("J3").Value = "belum"
$ws.Range("G3").Value = 261294
$ws.Range("J4").Value = "belum"
$ws.Range("G4").Value = 261294
$ws.Range("J5").Value = "belum"
$ws.Range("G5").Value = 261294
$ws.Range("J6").Value = "belum"
$ws.Range("G6").Value = 261294
$ws.Range("J7").Value = "belum"
$ws.Range("B7").Value = "fairytailfg@gmail.com"
$ws.Range("D7").Value = "perempuan"
$ws.Range("E7").Value = "Application1!"
$ws.Range("F7").Value = "Application1!"
$ws.Range("G7").Value = 261294
$ws.Range("I7").Value = "failed"
$ws.Range("J8").Value = "belum"
$ws.Range("B8").Value = "fairytailfg@gmail.com"
$ws.Range("C8").Value = 822
$ws.Range("D8").Value = "laki"
$ws.Range("E8").Value = "Application1!"
$ws.Range("F8").Value = "Application1!"
$ws.Range("G8").Value = 261294
$ws.Range("I8").Value = "failed"
$ws.Range("J9").Value = "belum"
$ws.Range("B9").Value = "fairytailfg@gmail.com"
$ws.Range("C9").Value = 82297476950
$ws.Range("D9").Value = "perempuan"
$ws.Range("E9").Value = "Application1!"
$ws.Range("F9").Value = "Application1!"
$ws.Range("G9").Value = 261294
$ws.Range("I9").Value = "failed"
$ws.Range("J10").Value = "belum"
$ws.Range("B10").Value = "fairytailfg@gmail.com"
$ws.Range("C10").Value = 82297476950
$ws.Range("D10").Value = "laki"
$ws.Range("G10").Value = 261294
$ws.Range("I10").Value = "failed"
$ws.Range("J11").Value = "belum"
$ws.Range("B11").Value = "fairytailfg@gmail.com"
$ws.Range("C11").Value = 82297476950
$ws.Range("D11").Value = "perempuan"
$ws.Range("E11").Value = "Application1!"
$ws.Range("G11").Value = 261294
$ws.Range("I11").Value = "failed"
$ws.Range("J12").Value = "belum"
$ws.Range("B12").Value = "fairytailfg@gmail.com"
$ws.Range("C12").Value = 82297476950
$ws.Range("D12").Value = "laki"
$ws.Range("E12").Value = "Application1!"
$ws.Range("F12").Value = "Application1!"
$ws.Range("G12").Value = 261294
$ws.Range("J13").Value = "sudah"
$ws.Range("B13").Value = "fairytailfg@gmail.com"
$ws.Range("D13").Value = "laki"
$ws.Range("E13").Value = "Application1!"
$ws.Range("F13").Value = "Application1!"
$ws.Range("G13").Value = 261294
$ws.Range("I13").Value = "failed"
$ws.Range("J14").Value = "sudah"
$ws.Range("B14").Value = "fairytailfg@gmail.com"
$ws.Range("D14").Value = "perempuan"
$ws.Range("E14").Value = "Application1!"
$ws.Range("F14").Value = "Application1!"
$ws.Range("G14").Value = 261294
$ws.Range("I14").Value = "failed"
$ws.Range("H13").Value = "noPolEmpty"
$ws.Range("H14").Value = "noPolNotFound"
$ws.Range("A7").Value = "Aing Macan"
$ws.Range("A8").Value = "Aing Singa"
$ws.Range("A9").Value = "Aing Gajah"
$ws.Range("H7").Value = "phoneNoEmpty"
$ws.Range("H8").Value = "phoneNoShort"
$ws.Range("H9").Value = "phoneNoExist"
$ws.Range("H10").Value = "passInvalid"
$ws.Range("H11").Value = "confPassNotMatch"
$ws.Range("A10").Value = "Tikus"
$ws.Range("A11").Value = "Jerapah"
$ws.Range("F11").Value = "Application2!"
$ws.Range("A12").Value = "Yahoo"
$ws.Range("I12").Value = "succeed"
$ws.Range("K14").Value = "w1316pk"
$ws.Range("E10").Value = "application"
$ws.Range("F10").Value = "application"

# ---------------------------------------------------------------------------
# Hyperlinks for the new B-column e-mail cells (order matches the source
# workbook's relationship ids: B14, B13, B7, B8, B9, B10, B11, B12)
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("B14"), "mailto:fairytailfg@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B13"), "mailto:fairytailfg@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B7"), "mailto:fairytailfg@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B8"), "mailto:fairytailfg@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B9"), "mailto:fairytailfg@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B10"), "mailto:fairytailfg@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B11"), "mailto:fairytailfg@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B12"), "mailto:fairytailfg@gmail.com")

# Hyperlinks.Add stamps its own ad-hoc "Hyperlink" style xf; re-apply the
# workbook's existing named "Hyperlink" style so every linked cell in B
# references the same style index as the original B2/B3/B6 cells.
$ws.Range("B7:B14").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# Column widths: new column A, and widened column H
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 10.666666666666666
$ws.Columns.Item(8).ColumnWidth = 16.666666666666668

# ---------------------------------------------------------------------------
# Sheet view: drop the frozen "topLeftCell" and move the selection to H9
# ---------------------------------------------------------------------------
$ws.Range("H9").Select()
